$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 25000
$ws.Range("J3").Value = 25000
$ws.Range("L3").Value = 25000
$ws.Range("N3").Value = -25228
$ws.Range("H62").Value = 835506.3
$ws.Range("I62").Value = 1033958.2
$ws.Range("J62").Value = 174000
$ws.Range("K62").Value = 1033958.2
$ws.Range("L62").Value = 174000
$ws.Range("M62").Value = -1033334.2
$ws.Range("N62").Value = -175248
$ws.Range("H65").Value = 835506.3
$ws.Range("I65").Value = 1033958.2
$ws.Range("J65").Value = 174000
$ws.Range("K65").Value = 5169791
$ws.Range("L65").Value = 870000
$ws.Range("M65").Value = -5166671
$ws.Range("N65").Value = -876240
$ws.Range("H86").Value = 9140007
$ws.Range("I86").Value = 3586.25
$ws.Range("K86").Value = 3586.25
$ws.Range("M86").Value = -2463.25
$ws.Range("H88").Value = 2691.9443
$ws.Range("J88").Value = 2865.0833
$ws.Range("L88").Value = 2865.0833
$ws.Range("N88").Value = -3677.0833
$ws.Range("H89").Value = 9140007
$ws.Range("I89").Value = 3586.25
$ws.Range("K89").Value = 17931.25
$ws.Range("M89").Value = -12315.25
$ws.Range("H91").Value = 2691.9443
$ws.Range("J91").Value = 2865.0833
$ws.Range("L91").Value = 2865.0833
$ws.Range("N91").Value = -5673.0833
$ws.Range("H93").Value = 36000
$ws.Range("J93").Value = 36000
$ws.Range("L93").Value = 36000
$ws.Range("N93").Value = -40992
$ws.Range("H100").Value = 9190.286
$ws.Range("I100").Value = 10055.333
$ws.Range("K100").Value = 10055.333
$ws.Range("M100").Value = -9514.333000000001
$ws.Range("H102").Value = 25000
$ws.Range("J102").Value = 25000
$ws.Range("L102").Value = 25000
$ws.Range("N102").Value = -31490
$ws.Range("H137").Value = 19610394
$ws.Range("I137").Value = 90910184
$ws.Range("K137").Value = 272730552
$ws.Range("M137").Value = -272728002
$ws.Range("H138").Value = 1963
$ws.Range("I138").Value = 1498.25
$ws.Range("K138").Value = 4494.75
$ws.Range("M138").Value = 645.25

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2279
$ws.Range("I32").Value = 2370.9375
$ws.Range("J32").Value = 808
$ws.Range("K32").Value = 2370.9375
$ws.Range("L32").Value = 808
$ws.Range("M32").Value = -2083.9375
$ws.Range("N32").Value = -1382
$ws.Range("H61").Value = 2517.6
$ws.Range("I61").Value = 2439.5789
$ws.Range("J61").Value = 4000
$ws.Range("K61").Value = 2439.5789
$ws.Range("L61").Value = 4000
$ws.Range("M61").Value = -2227.5789
$ws.Range("N61").Value = -4424
$ws.Range("H132").Value = 1534.0769
$ws.Range("I132").Value = 1515.44
$ws.Range("K132").Value = 4546.32
$ws.Range("M132").Value = -2016.32
$ws.Range("H136").Value = 2517.6
$ws.Range("I136").Value = 2439.5789
$ws.Range("J136").Value = 4000
$ws.Range("K136").Value = 7318.736699999999
$ws.Range("L136").Value = 12000
$ws.Range("M136").Value = -4768.736699999999
$ws.Range("N136").Value = -17100

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H54").Value = 14316.667
$ws.Range("I54").Value = 7475
$ws.Range("K54").Value = 7475
$ws.Range("M54").Value = -6991
$ws.Range("H105").Value = 1876.5238
$ws.Range("I105").Value = 1688.5333
$ws.Range("K105").Value = 1688.5333
$ws.Range("M105").Value = 58.46669999999995
$ws.Range("H134").Value = 1549.5
$ws.Range("I134").Value = 1099
$ws.Range("J134").Value = 2000
$ws.Range("K134").Value = 3297
$ws.Range("L134").Value = 6000
$ws.Range("M134").Value = -762
$ws.Range("N134").Value = -11070

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1899.8
$ws.Range("I16").Value = 1874.75
$ws.Range("K16").Value = 1874.75
$ws.Range("M16").Value = -1587.75
$ws.Range("H99").Value = 2487.375
$ws.Range("I99").Value = 2342.7144
$ws.Range("K99").Value = 2342.7144
$ws.Range("M99").Value = -844.7143999999998
$ws.Range("H113").Value = 1899.8
$ws.Range("I113").Value = 1874.75
$ws.Range("K113").Value = 1874.75
$ws.Range("M113").Value = 295.25
$ws.Range("H126").Value = 2487.375
$ws.Range("I126").Value = 2342.7144
$ws.Range("K126").Value = 7028.1432
$ws.Range("M126").Value = -4558.1432
$ws.Range("H141").Value = 103005.36
$ws.Range("J141").Value = 103005.36
$ws.Range("L141").Value = 103005.36
$ws.Range("N141").Value = -113365.36

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 2750007.2
$ws.Range("I7").Value = 3666671.2
$ws.Range("K7").Value = 11000013.6
$ws.Range("M7").Value = -10999901.6
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("M64").ClearContents()
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("M67").ClearContents()
$ws.Range("N67").ClearContents()
$ws.Range("H80").Value = 4144
$ws.Range("I80").Value = 502
$ws.Range("J80").Value = 6875.5
$ws.Range("K80").Value = 1506
$ws.Range("L80").Value = 20626.5
$ws.Range("M80").Value = -570
$ws.Range("N80").Value = -22498.5
$ws.Range("H83").Value = 4144
$ws.Range("I83").Value = 502
$ws.Range("J83").Value = 6875.5
$ws.Range("K83").Value = 4518
$ws.Range("L83").Value = 61879.5
$ws.Range("M83").Value = 162
$ws.Range("N83").Value = -71239.5
$ws.Range("H92").Value = 326.77777
$ws.Range("I92").Value = 361.6
$ws.Range("J92").Value = 283.25
$ws.Range("K92").Value = 1084.8
$ws.Range("L92").Value = 849.75
$ws.Range("M92").Value = 163.1999999999998
$ws.Range("N92").Value = -3345.75
$ws.Range("H133").Value = 4907.273
$ws.Range("I133").Value = 4220
$ws.Range("K133").Value = 12660
$ws.Range("M133").Value = -7600

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 25000
$ws.Range("J20").Value = 25000
$ws.Range("L20").Value = 25000
$ws.Range("N20").Value = -25490
$ws.Range("H43").Value = 3428.5715
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()
$ws.Range("H80").Value = 3793.3142
$ws.Range("I80").Value = 3626.375
$ws.Range("K80").Value = 3626.375
$ws.Range("M80").Value = -2628.375
$ws.Range("H83").Value = 3793.3142
$ws.Range("I83").Value = 3626.375
$ws.Range("K83").Value = 18131.875
$ws.Range("M83").Value = -13139.875
$ws.Range("H132").Value = 4659.2974
$ws.Range("I132").Value = 3843.9656
$ws.Range("K132").Value = 11531.8968
$ws.Range("M132").Value = -9001.8968

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 83333.336
$ws.Range("J46").Value = 83333.336
$ws.Range("L46").Value = 83333.336
$ws.Range("N46").Value = -83795.336
$ws.Range("H122").Value = 2407.75
$ws.Range("I122").Value = 2499.3
$ws.Range("J122").Value = 1950
$ws.Range("K122").Value = 7497.900000000001
$ws.Range("L122").Value = 5850
$ws.Range("M122").Value = -5047.900000000001
$ws.Range("N122").Value = -10750
$ws.Range("H126").Value = 1536
$ws.Range("I126").Value = 1432.8889
$ws.Range("K126").Value = 4298.6667
$ws.Range("M126").Value = -1828.6667
$ws.Range("H134").Value = 83333.336
$ws.Range("J134").Value = 83333.336
$ws.Range("L134").Value = 250000.008
$ws.Range("N134").Value = -255070.008
$ws.Range("H136").Value = 2516.6287
$ws.Range("I136").Value = 985.3214
$ws.Range("K136").Value = 2955.9642
$ws.Range("M136").Value = -405.9642000000003
